$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The post "「大きくなったウソ」" (row 781) was removed from the source data.
# Delete that entire row; Excel shifts every following row (782-829) up by
# one, which also renumbers their A/B/C cell references automatically and
# shrinks the used range from A1:C829 to A1:C828.
$ws.Rows.Item(781).Delete()
